$wb = $excel.ActiveWorkbook

$income = $wb.Worksheets.Item("NKE Income Annual")
$balance = $wb.Worksheets.Item("NKE Balance Annual")

# Move "NKE Income Annual" so it comes before "NKE Balance Annual"
$income.Move($balance)

# Re-fetch the sheet by name since object handles can be position-bound.
$balance = $wb.Worksheets.Item("NKE Balance Annual")

# Activate the Balance sheet (now second tab) and set its selection to A3
$balance.Activate()
$balance.Range("A3").Select()
